$d = $word.ActiveDocument

# The two paragraph-pairs (PT + italic EN) under "Objetivos" and
# "Programa resumido" are being swapped with each other. Each of the
# four distinct strings occurs exactly once in the document. Find/Replace
# only rewrites text in place (it cannot relocate a paragraph's content
# to a different paragraph by itself), so the swap is done via temporary
# markers: first park each paragraph's current text behind a unique
# marker, then fill each marker with the text that belongs there after
# the swap.

# Step 1: park the current "Objetivos" texts (long PT/EN) behind markers.
$d.Content.Find.Execute(
    "Proporcionar ao estudante conhecimentos práticos nos processos tecnológicos de preparação de bebidas fermentadas e destiladas. Proporcionar ao estudante conhecimentos práticos nos processos tecnológicos de preparação de bebidas fermentadas e destiladas.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "##TEMP_OBJETIVOS_PT##", 2) | Out-Null

$d.Content.Find.Execute(
    "Provide the student with practical knowledge in the technological processes of preparing fermented and distilled beverages. Provide the student with practical knowledge in the technological processes of preparing fermented and distilled beverages.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "##TEMP_OBJETIVOS_EN##", 2) | Out-Null

# Step 2: park the current "Programa resumido" texts (short PT/EN)
# behind markers.
$d.Content.Find.Execute(
    "Elaboração prática de cerveja, aguardente, licores e iogurtes.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "##TEMP_RESUMIDO_PT##", 2) | Out-Null

$d.Content.Find.Execute(
    "Practical classes for preparing beer, cachaça and yogurts.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "##TEMP_RESUMIDO_EN##", 2) | Out-Null

# Step 3: the marker sitting in the "Objetivos" paragraphs gets filled
# with the short "Programa resumido" text (the swap target), and vice
# versa for the marker sitting in the "Programa resumido" paragraphs.
$d.Content.Find.Execute(
    "##TEMP_OBJETIVOS_PT##", $true, $false, $false, $false, $false, $true,
    1, $false,
    "Elaboração prática de cerveja, aguardente, licores e iogurtes.", 2) |
    Out-Null

$d.Content.Find.Execute(
    "##TEMP_OBJETIVOS_EN##", $true, $false, $false, $false, $false, $true,
    1, $false,
    "Practical classes for preparing beer, cachaça and yogurts.", 2) |
    Out-Null

$d.Content.Find.Execute(
    "##TEMP_RESUMIDO_PT##", $true, $false, $false, $false, $false, $true,
    1, $false,
    "Proporcionar ao estudante conhecimentos práticos nos processos tecnológicos de preparação de bebidas fermentadas e destiladas. Proporcionar ao estudante conhecimentos práticos nos processos tecnológicos de preparação de bebidas fermentadas e destiladas.",
    2) | Out-Null

$d.Content.Find.Execute(
    "##TEMP_RESUMIDO_EN##", $true, $false, $false, $false, $false, $true,
    1, $false,
    "Provide the student with practical knowledge in the technological processes of preparing fermented and distilled beverages. Provide the student with practical knowledge in the technological processes of preparing fermented and distilled beverages.",
    2) | Out-Null
